$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.649.70"
$ws.Range("E2").Value = "  +2.71%  "

# Row 3
$ws.Range("D3").Value = "1.862.35"
$ws.Range("E3").Value = "  +2.01%  "

# Row 4
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'245.72"
$ws.Range("E5").Value = "  +2.68%  "

# Row 6
$ws.Range("D6").Value = "'0.6996"
$ws.Range("E6").Value = "  +1.27%  "

# Row 7
$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").Value = "'0.07739"
$ws.Range("E8").Value = "  +1.64%  "

# Row 9
$ws.Range("D9").Value = "'0.3071"
$ws.Range("E9").Value = "  +1.90%  "

# Row 10
$ws.Range("D10").Value = "'23.66"
$ws.Range("E10").Value = "  +1.41%  "

# Row 11
$ws.Range("D11").Value = "'0.07774"
$ws.Range("E11").Value = "  +0.73%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.863.68"
$ws.Range("E12").Value = "  +1.93%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.163"
$ws.Range("E13").Value = "  +2.48%  "

# Row 14
$ws.Range("D14").Value = "'92.36"
$ws.Range("E14").Value = "  +2.76%  "

# Row 15
$ws.Range("D15").Value = "'0.6930"
$ws.Range("E15").Value = "  +3.31%  "

# Row 16
$ws.Range("D16").Value = "'6.569"
$ws.Range("E16").Value = "  +2.42%  "

# Row 17
$ws.Range("D17").Value = "29.633.48"
$ws.Range("E17").Value = "  +2.67%  "

# Row 18
$ws.Range("D18").Value = "'0.000008374"
$ws.Range("E18").Value = "  +1.20%  "

# Row 19
$ws.Range("D19").Value = "2.109.48"
$ws.Range("E19").Value = "  +1.86%  "

# Row 20
$ws.Range("D20").Value = "'242.01"
$ws.Range("E20").Value = "  -0.09%  "

# Row 21
$ws.Range("E21").Value = "  +1.38%  "

# Row 22
$ws.Range("D22").Value = "'0.9997"
$ws.Range("E22").Value = "  +0.02%  "

# Row 23
$ws.Range("D23").Value = "'7.621"
$ws.Range("E23").Value = "  +3.39%  "

# Row 24
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  +0.04%  "

# Row 25
$ws.Range("D25").Value = "'0.1511"
$ws.Range("E25").Value = "  +2.98%  "

# Row 26
$ws.Range("D26").Value = "'8.923"
$ws.Range("E26").Value = "  +2.47%  "

# Row 27
$ws.Range("D27").Value = "'159.58"
$ws.Range("E27").Value = "  -0.59%  "

# Row 28
$ws.Range("E28").Value = "  +1.10%  "

# Row 29
$ws.Range("D29").Value = "'1.542"
$ws.Range("E29").Value = "  +1.05%  "

# Row 30
$ws.Range("D30").Value = "'4.263"
$ws.Range("E30").Value = "  +2.00%  "

# Row 31
$ws.Range("D31").Value = "'4.199"
$ws.Range("E31").Value = "  +1.84%  "

# Row 32
$ws.Range("E32").Value = "  +0.12%  "

# Row 33
$ws.Range("D33").Value = "'0.05113"
$ws.Range("E33").Value = "  +0.37%  "

# Row 34
$ws.Range("D34").Value = "'0.7846"
$ws.Range("E34").Value = "  +5.54%  "

# Row 35
$ws.Range("E35").Value = "  +5.31%  "

# Row 36
$ws.Range("D36").Value = "'1.158"

# Row 37
$ws.Range("E37").Value = "  +0.18%  "

# Row 38
$ws.Range("D38").Value = "1.332.96"
$ws.Range("E38").Value = "  +11.54%  "

# Row 39
$ws.Range("D39").Value = "'0.01882"
$ws.Range("E39").Value = "  +2.84%  "

# Row 40
$ws.Range("D40").Value = "'2.736"
$ws.Range("E40").Value = "  +2.38%  "

# Row 41
$ws.Range("D41").Value = "'0.9670"
$ws.Range("E41").Value = "  +6.08%  "

# Row 42
$ws.Range("D42").Value = "'5.966"
$ws.Range("E42").Value = "  +14.58%  "

# Row 43
$ws.Range("D43").Value = "'106.55"
$ws.Range("E43").Value = "  -1.30%  "

# Row 44
$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  +0.08%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'9.773"
$ws.Range("E45").Value = "  +3.64%  "

# Row 46
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "2.009.23"
$ws.Range("E46").Value = "  +1.62%  "

# Row 47
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5216"
$ws.Range("E47").Value = "  +1.09%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'64.88"
$ws.Range("E48").Value = "  +4.37%  "

# Row 49
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.790"
$ws.Range("E49").Value = "  +3.88%  "

# Row 50
$ws.Range("B50").Value = "Aptos"
$ws.Range("C50").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D50").Value = "'7.005"
$ws.Range("E50").Value = "  +2.21%  "

# Row 51
$ws.Range("B51").Value = "TheSandbox"
$ws.Range("C51").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D51").Value = "'0.4206"
$ws.Range("E51").Value = "  +0.50%  "
